# week 5,6: fix typos
#
# 1) "An unreferenced object is an object that cannot be access from
#     anywhere else in the program" -> "... cannot be accessed from ..."
#    on slide 14 (Garbage Collection / Content Placeholder 2).
#
# 2) The cached display text of the auto-updating Date placeholder
#    fields (shown on the slide master, every slide layout and the
#    notes master) was refreshed from "9/26/2016" to "2/20/17" the
#    next time the deck was opened/saved in PowerPoint.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1) Fix the "cannot be access" -> "cannot be accessed" typo.
# ---------------------------------------------------------------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if (-not $shape.HasTextFrame) { continue }
        $tr = $shape.TextFrame.TextRange
        $full = $tr.Text
        $needle = "cannot be access from"
        $idx = $full.IndexOf($needle)
        if ($idx -ge 0) {
            # Replace just "access " (the word + trailing space) with
            # "accessed " so the run split lines up with "... cannot be "
            # | "accessed " | "from anywhere ...".
            $wordIdx = $full.IndexOf("access from")
            $sub = $tr.Characters($wordIdx + 1, 7)
            $sub.Text = "accessed "
        }
    }
}

# ---------------------------------------------------------------
# 2) Refresh the cached "9/26/2016" date-placeholder text to "2/20/17"
#    wherever it appears: slide master, every slide layout, and the
#    notes master.
# ---------------------------------------------------------------
function Update-DateText($shapes) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $shp = $shapes.Item($k)
        if (-not $shp.HasTextFrame) { continue }
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "9/26/2016") {
            $tr.Text = "2/20/17"
        }
    }
}

$master = $p.SlideMaster
Update-DateText $master.Shapes

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $cl = $layouts.Item($i)
    Update-DateText $cl.Shapes
}

$notesMaster = $p.NotesMaster
Update-DateText $notesMaster.Shapes
